# Re-run of demo-clin pvalues-by-histology stats ("rerun with NF1 intronic plp")
# Updates the tmb p-value (row 9, term = "tmb") on several histology sheets,
# and several rows (age/sex/race/ancestry/ethnicity/OS_years) on the
# "Neurofibroma plexiform" sheet.

$wb = $excel.ActiveWorkbook

# term = tmb (row 9) p-value updates, one per histology worksheet
$tmbUpdates = @{
    "Low-grade glioma"           = 0.295225120835944
    "Other tumor"                = 0.672416118643011
    "Medulloblastoma"            = 0.0114236813806622
    "Mixed neuronal-glial tumor" = 0.47293154793839
    "Ependymoma"                 = 0.233368156730013
    "Other high-grade glioma"    = 0.840678123585337
    "Craniopharyngioma"          = 0.434874810217698
    "DIPG or DMG"                = 0.762129815426868
}

foreach ($sheetName in $tmbUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B9").Value = $tmbUpdates[$sheetName]
}

# Neurofibroma plexiform sheet: several rows change (age, sex, race,
# ancestry, ethnicity all -> 0.589743589743591; sex row also updates)
$nf = $wb.Worksheets.Item("Neurofibroma plexiform")
$nf.Range("B2").Value = 0.225641025641026
$nf.Range("B4").Value = 0.589743589743591
$nf.Range("B5").Value = 0.589743589743591
$nf.Range("B6").Value = 0.476190476190476
$nf.Range("B7").Value = 0.456043956043956
$nf.Range("B8").Value = 0.291208791208791
